$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised weekly "Waargenomen" (observed) figures for several existing weeks ---
$ws.Range("G8").Value = 4303
$ws.Range("G11").Value = 2986
$ws.Range("G14").Value = 2726
$ws.Range("G20").Value = 2618
$ws.Range("G22").Value = 2672
$ws.Range("G24").Value = 2637
$ws.Range("G26").Value = 2843
$ws.Range("G27").Value = 2730
$ws.Range("G28").Value = 2679
$ws.Range("G29").Value = 2730
$ws.Range("G30").Value = 2705
$ws.Range("G31").Value = 2878
$ws.Range("G32").Value = 2977

# --- New week 41 row ---
$ws.Range("F33").Value = 41
$ws.Range("G33").Value = 2959
$ws.Range("H33").Value = 2807
$ws.Range("I33").Formula = "=G33-H33"

# --- Old totals row (was row 35) gets removed entirely ... ---
$ws.Range("F35:I35").Clear()

# --- ... and recreated two rows down (row 37), same formulas, refreshed values ---
$ws.Range("F37").Value = "Som week 11 tot en met 19"
$ws.Range("G37").Formula = "=SUM(G3:G28)"
$ws.Range("G37").NumberFormat = "0"
$ws.Range("H37").Formula = "=SUM(H3:H28)"
$ws.Range("H37").NumberFormat = "0"
$ws.Range("I37").Formula = "=SUM(I3:I28)"
$ws.Range("I37").NumberFormat = "0"

# --- Match the author's final selection in the saved workbook ---
$ws.Range("J14").Select()
